$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 3575.7778
$ws.Range("I28").Value = 363.83334
$ws.Range("J28").Value = 9999.666999999999
$ws.Range("K28").Value = 363.83334
$ws.Range("L28").Value = 9999.666999999999
$ws.Range("M28").Value = 121.16666
$ws.Range("N28").Value = -10969.667
$ws.Range("H70").Value = 64121.117
$ws.Range("J70").Value = 96742.73
$ws.Range("L70").Value = 290228.19
$ws.Range("N70").Value = -290768.19
$ws.Range("H73").Value = 64121.117
$ws.Range("J73").Value = 96742.73
$ws.Range("L73").Value = 290228.19
$ws.Range("N73").Value = -292100.19
$ws.Range("H88").Value = 4358.125
$ws.Range("J88").Value = 4218.1816
$ws.Range("L88").Value = 4218.1816
$ws.Range("N88").Value = -5030.1816
$ws.Range("H91").Value = 4358.125
$ws.Range("J91").Value = 4218.1816
$ws.Range("L91").Value = 4218.1816
$ws.Range("N91").Value = -7026.1816
$ws.Range("H113").Value = 4128.4287
$ws.Range("I113").Value = 3499.6667
$ws.Range("J113").Value = 4600
$ws.Range("K113").Value = 3499.6667
$ws.Range("L113").Value = 4600
$ws.Range("M113").Value = -245.6667000000002
$ws.Range("N113").Value = -11108
$ws.Range("H115").Value = 2789.4
$ws.Range("I115").Value = 2789.4
$ws.Range("K115").Value = 8368.200000000001
$ws.Range("M115").Value = -6801.200000000001
$ws.Range("H127").Value = 555
$ws.Range("I127").Value = 555
$ws.Range("K127").Value = 1665
$ws.Range("M127").Value = 3295
$ws.Range("H137").Value = 1820.84
$ws.Range("J137").Value = 2765.6365
$ws.Range("L137").Value = 8296.9095
$ws.Range("N137").Value = -13396.9095
$ws.Range("H138").Value = 3468.7544
$ws.Range("I138").Value = 1220.92
$ws.Range("J138").Value = 5224.875
$ws.Range("K138").Value = 3662.76
$ws.Range("L138").Value = 15674.625
$ws.Range("M138").Value = 1477.24
$ws.Range("N138").Value = -25954.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 25344
$ws.Range("J96").Value = 25344
$ws.Range("L96").Value = 25344
$ws.Range("N96").Value = -30836
$ws.Range("H102").Value = 2721.4644
$ws.Range("I102").Value = 2824.76
$ws.Range("K102").Value = 2824.76
$ws.Range("M102").Value = -1202.76

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 40699.8
$ws.Range("J76").Value = 40699.8
$ws.Range("L76").Value = 40699.8
$ws.Range("N76").Value = -41329.8
$ws.Range("H79").Value = 40699.8
$ws.Range("J79").Value = 40699.8
$ws.Range("L79").Value = 40699.8
$ws.Range("N79").Value = -42883.8
$ws.Range("H105").Value = 2416.5
$ws.Range("I105").Value = 2554.5715
$ws.Range("J105").Value = 1450
$ws.Range("K105").Value = 2554.5715
$ws.Range("L105").Value = 1450
$ws.Range("M105").Value = -807.5715
$ws.Range("N105").Value = -4944

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3663.9756
$ws.Range("I31").Value = 3647.5557
$ws.Range("K31").Value = 3647.5557
$ws.Range("M31").Value = -3352.5557
$ws.Range("H34").Value = 3663.9756
$ws.Range("I34").Value = 3647.5557
$ws.Range("K34").Value = 3647.5557
$ws.Range("M34").Value = -3445.5557
$ws.Range("H105").Value = 1906.5
$ws.Range("I105").Value = 1750.2858
$ws.Range("K105").Value = 1750.2858
$ws.Range("M105").Value = -3.285800000000108
$ws.Range("H134").Value = 4776.857
$ws.Range("I134").Value = 4571.2573
$ws.Range("J134").Value = 5804.857
$ws.Range("K134").Value = 13713.7719
$ws.Range("L134").Value = 17414.571
$ws.Range("M134").Value = -11178.7719
$ws.Range("N134").Value = -22484.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 77288130
$ws.Range("J4").Value = 810800
$ws.Range("L4").Value = 2432400
$ws.Range("N4").Value = -2432624
$ws.Range("H45").Value = 9900
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H75").Value = 2151.7693
$ws.Range("I75").Value = 2358.25
$ws.Range("J75").Value = 2060
$ws.Range("K75").Value = 7074.75
$ws.Range("L75").Value = 6180
$ws.Range("M75").Value = -6076.75
$ws.Range("N75").Value = -8176
$ws.Range("H76").Value = 16763
$ws.Range("I76").Value = 8526
$ws.Range("K76").Value = 25578
$ws.Range("M76").Value = -25195
$ws.Range("H78").Value = 2151.7693
$ws.Range("I78").Value = 2358.25
$ws.Range("J78").Value = 2060
$ws.Range("K78").Value = 21224.25
$ws.Range("L78").Value = 18540
$ws.Range("M78").Value = -16232.25
$ws.Range("N78").Value = -28524
$ws.Range("H79").Value = 16763
$ws.Range("I79").Value = 8526
$ws.Range("K79").Value = 25578
$ws.Range("M79").Value = -24252
$ws.Range("H81").Value = 3099.4
$ws.Range("I81").Value = 2833.3333
$ws.Range("J81").Value = 3498.5
$ws.Range("K81").Value = 8499.999899999999
$ws.Range("L81").Value = 10495.5
$ws.Range("M81").Value = -7376.999899999999
$ws.Range("N81").Value = -12741.5
$ws.Range("H84").Value = 3099.4
$ws.Range("I84").Value = 2833.3333
$ws.Range("J84").Value = 3498.5
$ws.Range("K84").Value = 25499.9997
$ws.Range("L84").Value = 31486.5
$ws.Range("M84").Value = -19883.9997
$ws.Range("N84").Value = -42718.5
$ws.Range("H92").Value = 351.66666
$ws.Range("J92").Value = 373.6
$ws.Range("L92").Value = 1120.8
$ws.Range("N92").Value = -3616.8
$ws.Range("H107").Value = 557
$ws.Range("J107").Value = 530.3
$ws.Range("L107").Value = 1590.9
$ws.Range("N107").Value = -5430.9
$ws.Range("H112").Value = 6446.6875
$ws.Range("I112").Value = 432
$ws.Range("J112").Value = 7068.8965
$ws.Range("K112").Value = 1296
$ws.Range("L112").Value = 21206.6895
$ws.Range("M112").Value = -188
$ws.Range("N112").Value = -23422.6895
$ws.Range("H119").Value = 7597.6
$ws.Range("J119").Value = 19997
$ws.Range("L119").Value = 59991
$ws.Range("N119").Value = -69667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1503.8125
$ws.Range("I122").Value = 1503.8125
$ws.Range("K122").Value = 4511.4375
$ws.Range("M122").Value = -2061.4375
$ws.Range("H132").Value = 4146.773
$ws.Range("I132").Value = 3980.1177
$ws.Range("K132").Value = 11940.3531
$ws.Range("M132").Value = -9410.3531

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2244.5625
$ws.Range("J22").Value = 2483.5715
$ws.Range("L22").Value = 2483.5715
$ws.Range("N22").Value = -3073.5715
$ws.Range("H27").Value = 2244.5625
$ws.Range("J27").Value = 2483.5715
$ws.Range("L27").Value = 2483.5715
$ws.Range("N27").Value = -2697.5715
$ws.Range("H64").Value = 11499.5
$ws.Range("J64").Value = 12999
$ws.Range("L64").Value = 12999
$ws.Range("N64").Value = -13449
$ws.Range("H67").Value = 11499.5
$ws.Range("J67").Value = 12999
$ws.Range("L67").Value = 12999
$ws.Range("N67").Value = -14559
$ws.Range("H100").Value = 5799.9
$ws.Range("I100").Value = 4999.5
$ws.Range("K100").Value = 4999.5
$ws.Range("M100").Value = -4458.5
$ws.Range("H136").Value = 2069.6667
$ws.Range("I136").Value = 1914.7
$ws.Range("K136").Value = 5744.1
$ws.Range("M136").Value = -3194.1
